$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A slightly (13.85546875 -> 14.85546875)
$ws.Columns.Item(1).ColumnWidth = 14

# New rows of data (rows 7-18), matching the pattern of existing rows
$data = @(
    @(42600.782094907408, "Named", 13012, 5465, 329, 51, 20, 71, 28, 0, 0, 0, 0),
    @(42600.804664351854, "Named", 12598, 5260, 325, 40, 18, 68, 30, 0, 0, 0, 0),
    @(42600.806817129633, "Named", 12109, 5266, 326, 30, 18, 62, 37, 0, 0, 0, 0),
    @(42600.823229166665, "Named", 10245, 4852, 301, 26, 17, 60, 39, 0, 0, 0, 0),
    @(42600.830243055556, "Named", 12770, 5197, 322, 29, 18, 61, 38, 0, 0, 0, 0),
    @(42600.84170138889,  "Named", 11644, 5469, 331, 30, 21, 58, 41, 0, 0, 0, 0),
    @(42600.861226851855, "Named", 10121, 4888, 297, 28, 19, 59, 40, 0, 0, 0, 0),
    @(42600.868344907409, "Named", 10340, 4698, 279, 26, 19, 57, 42, 0, 0, 0, 0),
    @(42600.878541666665, "Named", 10578, 5161, 306, 41, 22, 64, 34, 0, 0, 0, 0),
    @(42600.883969907409, "Named", 10429, 5145, 313, 40, 22, 64, 35, 0, 0, 0, 0),
    @(42600.885185185187, "Named", 11286, 5149, 315, 40, 22, 64, 35, 0, 0, 0, 0),
    @(42600.886354166665, "Named", 10306, 4791, 292, 39, 19, 67, 32, 0, 0, 0, 0)
)

$startRow = 7
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"
    for ($c = 2; $c -le 13; $c++) {
        $ws.Cells.Item($row, $c).Value = $rowData[$c - 1]
    }
}
